$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Rules")

# Change cell B11 value from "R40" to the (text) string "1"
$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
